$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Fix the "Instituto Nacional de Migracion" source name (A9)
# ---------------------------------------------------------------
$ws.Range("A9").Value = "Instituto Nacional de Migración (INM)"

# ---------------------------------------------------------------
# 2) Row 10 - Comision Nacional de Bancos y Seguro (CNBS)
# ---------------------------------------------------------------
$ws.Range("A10").Value = "Comision Nacional de Bancos y Seguro (CNBS)"
$ws.Range("C10").Value = "Trabajo"
$ws.Range("D10").Value = "Institución que por mandato constitucional tiene la responsabilidad de velar por la estabilidad y solvencia del sistema financiero y demás supervisados, su regulación, supervisión y control. Asimismo, vigilamos la transparencia y que se respeten los derechos de los usuarios financieros, así como coadyuvamos con el sistema de prevención y detección del lavado activos y financiamiento al terrorismo, y contribuimos a promover la educación e inclusión financiera, a fin de salvaguardar el interés público."

$ws.Range("E10").Value = "https://covid19honduras.org/?q=cnbs-22-3"
$ws.Hyperlinks.Add($ws.Range("E10"), "https://covid19honduras.org/?q=cnbs-22-3") | Out-Null
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null

$ws.Range("F10").Value = "Las instituciones por la CNBS que realizan operaciones de crédito, podran otorgar periodos de gracia a los deudores que sean afectados por la reduccion de sus flujos de efectivo los cuales se podran otorgar hasta el 30 de junio de 2020."

$ws.Range("G10").Value = "https://covid19honduras.org/?q=cnbs-22-3"
$ws.Hyperlinks.Add($ws.Range("G10"), "https://covid19honduras.org/?q=cnbs-22-3") | Out-Null
$ws.Range("G9").Copy() | Out-Null
$ws.Range("G10").PasteSpecial(-4122) | Out-Null

$ws.Range("H10").Value = "21/3/2020"
$ws.Range("I10").Value = "22/3/2020"
$ws.Range("J10").Value = "Honduras"

$ws.Rows.Item(10).RowHeight = 105

# ---------------------------------------------------------------
# 3) Row 11 - Secretaria de Trabajo y Seguridad Social
# ---------------------------------------------------------------
$ws.Range("A11").Value = "Secretaria de Trabajo y Seguridad Social"
$ws.Range("C11").Value = "Trabajo"

$ws.Range("G11").Value = "https://covid19honduras.org/?q=secretaria-de-trabajo"
$ws.Hyperlinks.Add($ws.Range("G11"), "https://covid19honduras.org/?q=secretaria-de-trabajo") | Out-Null
$ws.Range("G9").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null

$ws.Range("I11").Value = "26/3/2020"
$ws.Range("J11").Value = "Honduras"

$ws.Rows.Item(11).RowHeight = 45

# ---------------------------------------------------------------
# 4) Selection / view state
# ---------------------------------------------------------------
$ws.Range("G11").Select() | Out-Null

$excel.CutCopyMode = $false
